$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Name"
$ws.Range("C2").Value = "Andi"
$ws.Range("C3").Value = "Eka"

$ws.Range("F6").Select()
